$wb = $excel.ActiveWorkbook

# The existing "Croatia" sheet is the template for the new "Greece" market
# test-data sheet: select it, select all its cells (mirrors the
# Move-or-Copy "Create a copy" workflow), and copy it to the end of the
# workbook, producing sheet "Croatia (2)".
$croatia = $wb.Worksheets.Item("Croatia")
$croatia.Activate()
$croatia.Cells.Select()
$croatia.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))

# Rename the freshly created copy and point it at the Greece market data.
$greece = $wb.Worksheets.Item($wb.Worksheets.Count)
$greece.Name = "Greece"
$greece.Range("B2").Value = "Greece Market"
$greece.Range("B4").Value = "NGC-4119/T3164"

# Leave the new sheet active with B4 selected, matching the normal
# "just typed into B4" end state of the edit.
$greece.Activate()
$greece.Range("B4").Select()
